$d = $word.ActiveDocument

# Helper: return the 1-based index of the LAST paragraph whose text starts
# with $text (re-scanned fresh each call so it stays correct across edits).
function Find-ParaIndex($text) {
    $idx = 0
    $foundIdx = -1
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like "$text*") {
            $foundIdx = $idx
        }
    }
    return $foundIdx
}

# --- 1) New bullet item inserted at the top of the "To be added:" list ---
# "Slide to left to delete an event" used to be the first bullet; a new
# bullet ("To remember when  app is closed") is inserted above it, and
# (as happens in Word when a list item is created this way) it gets its
# own fresh list numbering instance rather than re-using numId 14.
$slideIdx = Find-ParaIndex "Slide to left to delete an event"
$d.Paragraphs.Item($slideIdx).Range.InsertParagraphBefore()
$newPara1 = $d.Paragraphs.Item($slideIdx)
$newPara1.Range.Text = "To remember when  app is closed"
$newList = $d.ListTemplates.Add()
$newPara1.Range.ListFormat.ApplyListTemplateWithLevel($newList)

# --- 2) New bullet item appended after "Sorting by importance" ---
# "Notes, to note important things" joins the existing list, keeping the
# same numId (14) as its neighbours.
$sortIdx = Find-ParaIndex "Sorting by importance"
$d.Paragraphs.Item($sortIdx).Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($sortIdx + 1)
$newPara2.Range.Text = "Notes, to note important things"
